$wb = $excel.ActiveWorkbook

# --- Update "Logs" sheet: append new row 26 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A26").Value = "Ruilen van product"
$logs.Range("B26").Value = "mailmind.test@zohomail.eu"
$logs.Range("C26").Value = "Kan ik ruilen voor een andere maat?`nSent using {0}"
$logs.Range("D26").Value = "Retour / Terugbetaling"
$logs.Range("E26").Value = "Beste klant,`nBedankt voor je e-mail. Ja, het is mogelijk om te ruilen voor een andere maat. Graag ontvangen wij je ordernummer en de maat die je wilt ruilen, zodat we dit verder kunnen afhandelen. Aarzel niet om contact met ons op te nemen als je nog verdere vragen hebt.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Range("F26").Value = "2025-06-24 21:04:41"
$logs.Range("G26").Value = "Ja"

# --- Extend the conditional formatting ranges to cover the new row ---
$dFc = $logs.Range("D2:D25").FormatConditions
$dFc.Item(1).ModifyAppliesToRange($logs.Range("D2:D26"))

$gFc = $logs.Range("G2:G25").FormatConditions
$gFc.Item(1).ModifyAppliesToRange($logs.Range("G2:G26"))

# --- Update "Dashboard" sheet: bump the "Retour / Terugbetaling" count ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 9
